$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 used to hold the "TrendingUrlTemplate" setting; replace it with a new
# "YahooFinanceUrl" setting whose value is a live hyperlink.
$ws.Range("A9").Value = "YahooFinanceUrl"
$ws.Range("B9").Value = "https://finance.yahoo.com/"
$ws.Hyperlinks.Add($ws.Range("B9"), "https://finance.yahoo.com/")
# Re-apply the shared "Hyperlink" cell style so B9 matches the look of the
# other hyperlink cell (B10) instead of getting a brand new style entry.
$ws.Range("B9").Style = $ws.Range("B10").Style

# Move the selection to A9 (matches the saved cursor position in the file).
$ws.Range("A9").Select()
